$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.533.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.384.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.73%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '592.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.52%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.135'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.421'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.971.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.94%  '
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.71'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.584.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.398.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '449.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +13.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.32%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.519.53'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.54%  '
$ws.Range('E25').Value = '  +4.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.522'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.56'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.60'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.754.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.76%  '
$ws.Range('E44').Value = '  +3.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.69'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0692'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '343.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.38%  '
$ws.Range('E51').Value = '  +5.83%  '
